# Generate Report for Handoff
#
# A new handoff was produced for the "b.md" file (in both the zh-cn and
# de-de localization sheets). This updates:
#   - Status:              "Handed back: in sync with en-US" -> "Ready for handoff"
#   - Latest Handoff File:  a.6631f68b...<lang>.xlf           -> b.63290e5768f688058c7b37413b0a5c26c308f864.<lang>.xlf
#   - Latest Handoff Datetime updated to the new handoff timestamp
#
# The Overview sheet's summary status cells for "b.md" are updated the same way.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: b.md row (row 3) status columns ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"

# --- zh-cn sheet: b.md row (row 3) ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B3").Value = "Ready for handoff"
$zhcn.Range("C3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("D3").Value = "2016-01-28 04:02:09"

# --- de-de sheet: b.md row (row 3) ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B3").Value = "Ready for handoff"
$dede.Range("C3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("D3").Value = "2016-01-28 04:02:19"
